# Updated cryptos list prices/volumes (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Price (D) / Volume(1h) (E) text. $null means "leave D unchanged".
$updates = @{
    2 = @("30.937.14", "  +0.40%  ")
    3 = @("1.932.33", "  -0.09%  ")
    4 = @("0.9990", "  -0.21%  ")
    5 = @("242.24", "  -0.28%  ")
    6 = @("0.9990", "  -0.27%  ")
    7 = @("0.4791", "  -2.09%  ")
    8 = @("0.2885", "  -2.42%  ")
    9 = @("0.06790", "  -1.09%  ")
    10 = @("19.73", "  +2.49%  ")
    11 = @("104.54", "  -1.10%  ")
    12 = @("0.07816", "  +0.40%  ")
    13 = @("1.922.47", "  -0.73%  ")
    14 = @("5.291", "  -0.73%  ")
    15 = @("0.6833", "  -2.40%  ")
    16 = @("295.70", "  +8.23%  ")
    17 = @("30.934.27", "  +0.37%  ")
    18 = @("2.195.14", "  +0.06%  ")
    19 = @("0.000007600", "  -1.27%  ")
    20 = @("1.0000", "  -0.09%  ")
    21 = @("12.91", "  -0.99%  ")
    22 = @("5.525", "  -1.73%  ")
    23 = @("0.9992", "  -0.19%  ")
    24 = @("6.402", "  -1.71%  ")
    25 = @("9.558", "  -2.73%  ")
    26 = @("168.21", "  +2.11%  ")
    27 = @("19.83", "  +1.40%  ")
    28 = @("2.115", "  -2.20%  ")
    29 = @("1.392", "  +0.24%  ")
    30 = @("0.1015", "  -2.04%  ")
    31 = @("4.631", "  +1.57%  ")
    32 = @("1.532", "  -1.16%  ")
    33 = @("4.354", "  -0.72%  ")
    34 = @("0.04833", "  -0.99%  ")
    35 = @("0.7393", "  -2.24%  ")
    36 = @("1.129", "  -1.50%  ")
    37 = @("2.725", "  +0.07%  ")
    38 = @("0.01955", "  -2.20%  ")
    39 = @("2.634", "  -1.02%  ")
    40 = @("6.498", "  +0.37%  ")
    41 = @("76.99", "  -3.31%  ")
    42 = @("2.036", "  -1.35%  ")
    43 = @("0.8728", "  -1.76%  ")
    44 = @($null, "  -1.74%  ")
    45 = @("106.30", "  -1.77%  ")
    46 = @("0.9990", "  -0.26%  ")
    47 = @("1.026.33", "  +4.37%  ")
    48 = @("7.553", "  -4.12%  ")
    49 = @("0.1211", "  -2.65%  ")
    50 = @("9.091", "  -1.05%  ")
    51 = @("35.15", "  -2.84%  ")
}

foreach ($row in $updates.Keys) {
    $newPrice = $updates[$row][0]
    $newVolume = $updates[$row][1]

    if ($newPrice -ne $null) {
        # Force text storage so strings like "0.9990" or "1.0000" keep their
        # exact digits instead of being coerced into numbers.
        $priceCell = $ws.Range("D$row")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $newPrice
        $priceCell.ClearFormats()
    }

    $volumeCell = $ws.Range("E$row")
    $volumeCell.NumberFormat = "@"
    $volumeCell.Value = $newVolume
    $volumeCell.ClearFormats()
}
